$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 309, pushing existing rows 309:412 down to 310:413
$ws.Rows(309).Insert()

# Populate the newly inserted row 309 with the new record's data
$ws.Cells.Item(309, 1).Value = 10
$ws.Cells.Item(309, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(309, 3).Value = "La Araucanía"
$ws.Cells.Item(309, 4).Value = 44809
$ws.Cells.Item(309, 5).Value = 9
$ws.Cells.Item(309, 6).Value = "Fruta"
$ws.Cells.Item(309, 7).Value = 100108
$ws.Cells.Item(309, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(309, 9).Value = 100108002
$ws.Cells.Item(309, 10).Value = "Mango"
$ws.Cells.Item(309, 11).Value = "Sin especificar"
$ws.Cells.Item(309, 12).Value = "Primera"
$ws.Cells.Item(309, 13).Value = 350
$ws.Cells.Item(309, 14).Value = 10000
$ws.Cells.Item(309, 15).Value = 10000
$ws.Cells.Item(309, 16).Value = 10000
$ws.Cells.Item(309, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item(309, 18).Value = "Brasil"
$ws.Cells.Item(309, 19).Value = 2500
$ws.Cells.Item(309, 20).Value = 4
